$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F5:F11 -> 2, F12:F18 -> 3, F19:F25 -> 4, F26:F32 -> 5, F33:F39 -> 6, F40:F46 -> 7
$ws.Range("F5:F11").Value = 2
$ws.Range("F12:F18").Value = 3
$ws.Range("F19:F25").Value = 4
$ws.Range("F26:F32").Value = 5
$ws.Range("F33:F39").Value = 6
$ws.Range("F40:F46").Value = 7

# Update sheet view
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("F49").Select()
